# Re-sort the Local Authority -> Combined Authority lookup table.
# Sort by CAUTH23NM (col D) then LAD23NM (col B), matching the new
# LSIP / LEP / MCA naming and latest local authority list used
# elsewhere in the dashboard.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Local_Authority_District_to_Com")

$sortRange = $ws.Range("A2:E54")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D54"))
$ws.Sort.SortFields.Add($ws.Range("B2:B54"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Leave the selection where the editor left off after the sort.
$ws.Range("F15").Select()
